$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''69.841.00'
$ws.Range('E2').Value = '''  -1.80%  '
$ws.Range('D3').Value = '''3.789.27'
$ws.Range('E3').Value = '''  +2.87%  '
$ws.Range('E4').Value = '''  +0.02%  '
$ws.Range('D5').Value = '''622.14'
$ws.Range('E5').Value = '''  +3.99%  '
$ws.Range('D6').Value = '''177.84'
$ws.Range('E6').Value = '''  -3.62%  '
$ws.Range('D7').Value = '''3.787.21'
$ws.Range('E7').Value = '''  +2.88%  '
$ws.Range('E8').Value = '''  +0.08%  '
$ws.Range('D9').Value = '''0.535'
$ws.Range('E9').Value = '''  -0.26%  '
$ws.Range('D10').Value = '''0.170'
$ws.Range('E10').Value = '''  +3.61%  '
$ws.Range('D11').Value = '''6.29'
$ws.Range('E11').Value = '''  -5.22%  '
$ws.Range('E12').Value = '''  -1.69%  '
$ws.Range('D13').Value = '''40.78'
$ws.Range('E13').Value = '''  +2.30%  '
$ws.Range('D14').Value = '''0.0000260'
$ws.Range('E14').Value = '''  +2.24%  '
$ws.Range('D15').Value = '''4.423.25'
$ws.Range('E15').Value = '''  +2.84%  '
$ws.Range('D16').Value = '''3.791.84'
$ws.Range('E16').Value = '''  +2.79%  '
$ws.Range('D17').Value = '''69.937.32'
$ws.Range('E17').Value = '''  -1.78%  '
$ws.Range('E18').Value = '''  +0.67%  '
$ws.Range('D19').Value = '''7.59'
$ws.Range('E19').Value = '''  +1.07%  '
$ws.Range('D20').Value = '''16.84'
$ws.Range('E20').Value = '''  -0.77%  '
$ws.Range('D21').Value = '''507.18'
$ws.Range('E21').Value = '''  -1.69%  '
$ws.Range('D22').Value = '''9.65'
$ws.Range('E22').Value = '''  +4.54%  '
$ws.Range('D23').Value = '''0.727'
$ws.Range('E23').Value = '''  -2.31%  '
$ws.Range('D24').Value = '''2.54'
$ws.Range('E24').Value = '''  +4.48%  '
$ws.Range('D25').Value = '''87.16'
$ws.Range('E25').Value = '''  -0.29%  '
$ws.Range('D26').Value = '''13.22'
$ws.Range('E26').Value = '''  -1.89%  '
$ws.Range('D27').Value = '''11.16'
$ws.Range('E27').Value = '''  +1.16%  '
$ws.Range('D28').Value = '''0.0000139'
$ws.Range('E28').Value = '''  +25.55%  '
$ws.Range('E29').Value = '''  +0.11%  '
$ws.Range('D30').Value = '''2.50'
$ws.Range('E30').Value = '''  -0.92%  '
$ws.Range('D31').Value = '''2.91'
$ws.Range('E31').Value = '''  +5.08%  '
$ws.Range('D32').Value = '''7.79'
$ws.Range('E32').Value = '''  -4.05%  '
$ws.Range('D33').Value = '''31.30'
$ws.Range('E33').Value = '''  -1.22%  '
$ws.Range('E34').Value = '''  -1.02%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '''  +0.00%  '
$ws.Range('E36').Value = '''  +5.95%  '
$ws.Range('D37').Value = '''6.19'
$ws.Range('E37').Value = '''  +1.21%  '
$ws.Range('D38').Value = '''0.336'
$ws.Range('E38').Value = '''  -1.92%  '
$ws.Range('E39').Value = '''  +4.08%  '
$ws.Range('E40').Value = '''  -2.28%  '
$ws.Range('D41').Value = '''50.90'
$ws.Range('E41').Value = '''  -0.28%  '
$ws.Range('D42').Value = '''45.70'
$ws.Range('E42').Value = '''  +0.91%  '
$ws.Range('D43').Value = '''424.39'
$ws.Range('E43').Value = '''  +3.66%  '
$ws.Range('D44').Value = '''8.73'
$ws.Range('E44').Value = '''  -1.22%  '
$ws.Range('D45').Value = '''3.040.01'
$ws.Range('E45').Value = '''  -4.27%  '
$ws.Range('E46').Value = '''  +1.04%  '
$ws.Range('D47').Value = '''0.0363'
$ws.Range('E47').Value = '''  -1.39%  '
$ws.Range('D48').Value = '''27.44'
$ws.Range('E48').Value = '''  -2.79%  '
$ws.Range('B49').Value = '''Monero'
$ws.Range('C49').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '''138.64'
$ws.Range('E49').Value = '''  +0.71%  '
$ws.Range('B50').Value = '''USDe'
$ws.Range('C50').Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '''  -0.05%  '
$ws.Range('D51').Value = '''2.48'
$ws.Range('E51').Value = '''  +1.03%  '
